$d = $word.ActiveDocument

# After the "LOB1235: Impactos e Adeqüação Ambiental (Requisito)" paragraph
# the document contains a blank paragraph, a paragraph with the text
# "Ver no Jupiter Salvar em pdf Salvar em docx", another blank paragraph
# and a blank page-break paragraph that all need to be removed, while the
# two blank paragraphs that follow them must stay untouched.
#
# Locate the "Ver no Jupiter..." paragraph by its text and delete the
# range that spans from the blank paragraph right before it through the
# blank page-break paragraph right after it.
$jupiterIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $jupiterIndex = $i
    }
}

if ($jupiterIndex -gt 0) {
    $startPara = $d.Paragraphs($jupiterIndex - 1)
    $endPara = $d.Paragraphs($jupiterIndex + 2)
    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $deleteRange.Delete()
}
